# v0.6.2p: Can add/change terrain in GameMapEditor
# Adds new terrain rows (Sidewalk, Sand, Water variants) to the Terrains
# sheet's lookup table, shifting the existing 151+ rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: make room -----------------------------------------------
# 7 new rows for terrain codes 132,133,141-145 (Sidewalk + Sand), inserted
# just above the old "151" (Grass, light) row.
$ws.Rows.Item(17).Resize(7).Insert()

# 5 new rows for terrain codes 181-185 (Water), inserted just above the
# old "201" (Brick, blue) row (which is now at row 43 after the insert
# above).
$ws.Rows.Item(43).Resize(5).Insert()

# --- Step 2: populate the new cells -----------------------------------
# Write the new text cells in the same order the strings were first
# introduced (Sand, then Sidewalk, then Water) so the shared-string table
# ends up in the same append order as the authored workbook.

$ws.Range("B19").Value = 141
$ws.Range("C19").Value = "Sand, light"

$ws.Range("B20").Value = 142
$ws.Range("C20").Value = "Sand, tan"

$ws.Range("B21").Value = 143
$ws.Range("C21").Value = "Sand, dark"

$ws.Range("B22").Value = 144
$ws.Range("C22").Value = "Sand, line left"

$ws.Range("B23").Value = 145
$ws.Range("C23").Value = "Sand, line up"

$ws.Range("B17").Value = 132
$ws.Range("C17").Value = "Sidewalk, smooth"

$ws.Range("B18").Value = 133
$ws.Range("C18").Value = "Sidewalk, cracked"

$ws.Range("B42").Value = 181
$ws.Range("C42").Value = "Water, rocks"

$ws.Range("B43").Value = 182
$ws.Range("C43").Value = "Water, dirt"

$ws.Range("B44").Value = 183
$ws.Range("C44").Value = "Water, shallow"

$ws.Range("B45").Value = 184
$ws.Range("C45").Value = "Water, medium"

$ws.Range("B46").Value = 185
$ws.Range("C46").Value = "Water, deep"

# --- Step 3: fix up number formatting on the new Water rows -----------
# Rows 42-46 use the same "quote-prefixed, right aligned" number style as
# the other section-leading rows (e.g. row 41/179, row 55/211). Writing a
# plain numeric .Value resets formatting, so stamp the style back on
# afterwards by copying formats from a cell that already has it.
$ws.Range("B41").Copy()
$ws.Range("B42:B46").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 4: selection / scroll position, matching the saved view -----
$ws.Range("E43").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
